# update code tinh luong cho Quyen
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tổng hợp lương")

$ws.Range("C3").Value = 7852571.428571429
$ws.Range("C4").Value = -839285.7142857146
$ws.Range("C6").Value = 915000
$ws.Range("C8").Value = 4736428.571428571
$ws.Range("C9").Value = 495000
$ws.Range("C11").Value = 9743000
$ws.Range("C12").Value = 15257190.47619048
$ws.Range("C13").Value = 11689428.57142857
$ws.Range("C14").Value = -1584000
$ws.Range("C15").Value = 48265333.33333333
